# Update "想去人数" (number of people interested) counts on two sheets:
# "展览" and "全部类型" - each count incremented by 1 (except the 2164 row,
# which increments by 2, matching the authoritative diff).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F9").Value = 6990
$ws1.Range("F10").Value = 81
$ws1.Range("F13").Value = 7887
$ws1.Range("F18").Value = 2347
$ws1.Range("F28").Value = 2166
$ws1.Range("F32").Value = 80
$ws1.Range("F36").Value = 1436

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F12").Value = 6990
$ws4.Range("F13").Value = 81
$ws4.Range("F16").Value = 7887
$ws4.Range("F21").Value = 2347
$ws4.Range("F33").Value = 2166
$ws4.Range("F37").Value = 80
$ws4.Range("F42").Value = 1436
